# Weekly fruit/vegetable price update: insert 3 new report rows for
# "Cuatro cascos" pepper right above the existing "Cuatro cascos verde"
# entry (old row 31), pushing every row below it down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 31; everything that was on row 31
# onward (old rows 31-100) shifts down to rows 34-103.
$ws.Rows("31:33").Insert()

# New row 31: "Cuatro cascos" / Primera
$ws.Range("A31").Value = 12
$ws.Range("B31").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 45246
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = 100112002
$ws.Range("G31").Value = "Pimiento"
$ws.Range("H31").Value = "Cuatro cascos"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 10
$ws.Range("K31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 30000
$ws.Range("N31").Value = "$/caja 18 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 1667
$ws.Range("Q31").Value = 18
$ws.Range("R31").Value = "Hortaliza"

# New row 32: "Cuatro cascos" / Segunda
$ws.Range("A32").Value = 12
$ws.Range("B32").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 45246
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 100112002
$ws.Range("G32").Value = "Pimiento"
$ws.Range("H32").Value = "Cuatro cascos"
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 15
$ws.Range("K32").Value = 28000
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = 28000
$ws.Range("N32").Value = "$/caja 18 kilos"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 1556
$ws.Range("Q32").Value = 18
$ws.Range("R32").Value = "Hortaliza"

# New row 33: "Cuatro cascos" / Tercera
$ws.Range("A33").Value = 12
$ws.Range("B33").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44467
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100112002
$ws.Range("G33").Value = "Pimiento"
$ws.Range("H33").Value = "Cuatro cascos"
$ws.Range("I33").Value = "Tercera"
$ws.Range("J33").Value = 20
$ws.Range("K33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = 25000
$ws.Range("N33").Value = "$/caja 18 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1389
$ws.Range("Q33").Value = 18
$ws.Range("R33").Value = "Hortaliza"
